# Refresh current market-price derived columns (H:N) across all Leve profit sheets.
# Values mirror the latest scheduled-runner pull; no structural changes otherwise.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1300.9546
$ws.Range("I19").Value = 1940.1428
$ws.Range("J19").Value = 1002.6667
$ws.Range("K19").Value = 1940.1428
$ws.Range("L19").Value = 1002.6667
$ws.Range("M19").Value = -1765.1428
$ws.Range("N19").Value = -1352.6667
$ws.Range("H98").Value = 1988.3334
$ws.Range("I98").Value = 2056.4285
$ws.Range("J98").Value = 1750
$ws.Range("K98").Value = 2056.4285
$ws.Range("L98").Value = 1750
$ws.Range("M98").Value = -558.4285
$ws.Range("N98").Value = -4746
$ws.Range("H103").Value = 6678261.5
$ws.Range("I103").Value = 20033534
$ws.Range("J103").Value = 625
$ws.Range("K103").Value = 60100602
$ws.Range("L103").Value = 1875
$ws.Range("M103").Value = -60100016
$ws.Range("N103").Value = -3047
$ws.Range("H122").Value = 1988.3334
$ws.Range("I122").Value = 2056.4285
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 6169.2855
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -3719.2855
$ws.Range("N122").Value = -10150

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9856.786
$ws.Range("I32").Value = 10187.184
$ws.Range("K32").Value = 10187.184
$ws.Range("M32").Value = -9900.183999999999
$ws.Range("H63").Value = 62502028
$ws.Range("I63").Value = 71430616
$ws.Range("J63").Value = 1950
$ws.Range("K63").Value = 71430616
$ws.Range("L63").Value = 1950
$ws.Range("M63").Value = -71429930
$ws.Range("N63").Value = -3322
$ws.Range("H66").Value = 62502028
$ws.Range("I66").Value = 71430616
$ws.Range("J66").Value = 1950
$ws.Range("K66").Value = 357153080
$ws.Range("L66").Value = 9750
$ws.Range("M66").Value = -357149648
$ws.Range("N66").Value = -16614
$ws.Range("H74").Value = 18522666
$ws.Range("I74").Value = 35716670
$ws.Range("K74").Value = 35716670
$ws.Range("M74").Value = -35715796
$ws.Range("H77").Value = 18522666
$ws.Range("I77").Value = 35716670
$ws.Range("K77").Value = 178583350
$ws.Range("M77").Value = -178578982
$ws.Range("H102").Value = 2550
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 2550
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 2550
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -5794

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1950
$ws.Range("I29").Value = 1600
$ws.Range("J29").Value = 3000
$ws.Range("K29").Value = 1600
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = -1311
$ws.Range("N29").Value = -3578
$ws.Range("H94").Value = 1042.5
$ws.Range("I94").Value = 999.7222
$ws.Range("J94").Value = 1138.75
$ws.Range("K94").Value = 999.7222
$ws.Range("L94").Value = 1138.75
$ws.Range("M94").Value = -548.7222
$ws.Range("N94").Value = -2040.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 88.72221999999999
$ws.Range("I7").Value = 72.7
$ws.Range("J7").Value = 108.75
$ws.Range("K7").Value = 72.7
$ws.Range("L7").Value = 108.75
$ws.Range("M7").Value = 40.3
$ws.Range("N7").Value = -334.75
$ws.Range("H31").Value = 13891102
$ws.Range("I31").Value = 2308.739
$ws.Range("J31").Value = 333333340
$ws.Range("K31").Value = 2308.739
$ws.Range("L31").Value = 333333340
$ws.Range("M31").Value = -2013.739
$ws.Range("N31").Value = -333333930
$ws.Range("H34").Value = 13891102
$ws.Range("I34").Value = 2308.739
$ws.Range("J34").Value = 333333340
$ws.Range("K34").Value = 2308.739
$ws.Range("L34").Value = 333333340
$ws.Range("M34").Value = -2106.739
$ws.Range("N34").Value = -333333744
$ws.Range("H94").Value = 2835.1724
$ws.Range("I94").Value = 1601.0769
$ws.Range("K94").Value = 1601.0769
$ws.Range("M94").Value = -1150.0769
$ws.Range("H107").Value = 676.9583
$ws.Range("I107").Value = 465.77274
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 465.77274
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 1454.22726
$ws.Range("N107").Value = -6840
$ws.Range("H122").Value = 1670.3077
$ws.Range("I122").Value = 1706.5714
$ws.Range("J122").Value = 1518
$ws.Range("K122").Value = 5119.7142
$ws.Range("L122").Value = 4554
$ws.Range("M122").Value = -2669.7142
$ws.Range("N122").Value = -9454
$ws.Range("H132").Value = 26318574
$ws.Range("I132").Value = 29414052
$ws.Range("J132").Value = 7014
$ws.Range("K132").Value = 88242156
$ws.Range("L132").Value = 21042
$ws.Range("M132").Value = -88239626
$ws.Range("N132").Value = -26102

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 74.625
$ws.Range("I14").Value = 74.625
$ws.Range("K14").Value = 223.875
$ws.Range("M14").Value = -50.875
$ws.Range("H80").Value = 2317.3333
$ws.Range("J80").Value = 2380.4
$ws.Range("L80").Value = 7141.200000000001
$ws.Range("N80").Value = -9013.200000000001
$ws.Range("H83").Value = 2317.3333
$ws.Range("J83").Value = 2380.4
$ws.Range("L83").Value = 21423.6
$ws.Range("N83").Value = -30783.6
$ws.Range("H112").Value = 6786.778
$ws.Range("J112").Value = 6500
$ws.Range("L112").Value = 19500
$ws.Range("N112").Value = -21716
$ws.Range("H113").Value = 985.65
$ws.Range("I113").Value = 480
$ws.Range("J113").Value = 1491.3
$ws.Range("K113").Value = 1440
$ws.Range("L113").Value = 4473.9
$ws.Range("M113").Value = 730
$ws.Range("N113").Value = -8813.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 4450
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 7900
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 7900
$ws.Range("M12").Value = -860
$ws.Range("N12").Value = -8180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7486.2666
$ws.Range("I40").Value = 7229.4
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 7229.4
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -7093.4
$ws.Range("N40").Value = -8272
$ws.Range("H132").Value = 9265785
$ws.Range("I132").Value = 4061.6333
$ws.Range("J132").Value = 20842940
$ws.Range("K132").Value = 12184.8999
$ws.Range("L132").Value = 62528820
$ws.Range("M132").Value = -9654.8999
$ws.Range("N132").Value = -62533880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 35000
$ws.Range("J19").Value = 35000
$ws.Range("L19").Value = 35000
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -35348
$ws.Range("H100").Value = 1497.5454
$ws.Range("I100").Value = 1985.3334
$ws.Range("J100").Value = 912.2
$ws.Range("K100").Value = 3970.6668
$ws.Range("L100").Value = 1824.4
$ws.Range("M100").Value = -3429.6668
$ws.Range("N100").Value = -2906.4
$ws.Range("H107").Value = 684.8261
$ws.Range("I107").Value = 799.64703
$ws.Range("J107").Value = 359.5
$ws.Range("K107").Value = 2398.94109
$ws.Range("L107").Value = 1078.5
$ws.Range("M107").Value = -478.9410899999998
$ws.Range("N107").Value = -4918.5
$ws.Range("H126").Value = 7927
$ws.Range("I126").Value = 3180.7144
$ws.Range("J126").Value = 19001.666
$ws.Range("K126").Value = 9542.143199999999
$ws.Range("L126").Value = 57004.99800000001
$ws.Range("M126").Value = -7072.143199999999
$ws.Range("N126").Value = -61944.99800000001
